$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 8 ("Technology Used" - currently "Google Drive") and Slide 9
# ("Technology Used" - currently "Emacs Org-Mode") swap their bullet content.
# ---------------------------------------------------------------------------

$slide8 = $p.Slides.Item(8)
$slide9 = $p.Slides.Item(9)

$body8 = $slide8.Shapes.Item(2).TextFrame.TextRange
$body9 = $slide9.Shapes.Item(2).TextFrame.TextRange

# --- Slide 8: "Google Drive" / "Realtime collaborative editing." ---------
#     becomes "Emacs Org-Mode" / "Writing documentation ... org-mode ... all."
$para8_1 = $body8.Paragraphs(1)
$para8_1.Runs(1).Text = "Emacs Org-Mode"

$para8_2 = $body8.Paragraphs(2)
# Drop any extra runs left over from the previous content, back to front so
# the indices of the runs we still need to edit don't shift underneath us.
$count8 = $para8_2.Runs().Count
for ($i = $count8; $i -ge 2; $i--) {
    $para8_2.Runs($i).Text = ""
}
$para8_2.Runs(1).Text = "Writing documentation and planning deadlines, org-mode does it all."

# Re-split that single run into three runs so the middle word ("org-mode")
# can carry its own Consolas formatting, matching the target markup.
$full8 = $para8_2.Runs(1).Text
$start8 = $full8.IndexOf("org-mode") + 1
$para8_2.Characters($start8, 8).Font.Name = "Consolas"

# --- Slide 9: "Emacs Org-Mode" / "Writing documentation ... org-mode ... all."
#     becomes "Google Drive" / "Realtime collaborative editing."
$para9_1 = $body9.Paragraphs(1)
$para9_1.Runs(1).Text = "Google Drive"

$para9_2 = $body9.Paragraphs(2)
$count9 = $para9_2.Runs().Count
for ($i = $count9; $i -ge 2; $i--) {
    $para9_2.Runs($i).Text = ""
}
$para9_2.Runs(1).Text = "Realtime collaborative editing."

# ---------------------------------------------------------------------------
# The two notes pages ("g13a42880235_2_34" for slide 8, "g13a42880235_2_41"
# for slide 9) swap their Google-Shape id suffixes along with the content.
# ---------------------------------------------------------------------------

$notes8 = $slide8.NotesPage
$notes8.Shapes.Item(1).Name = "Google Shape;173;g13a42880235_2_41:notes"
$notes8.Shapes.Item(2).Name = "Google Shape;174;g13a42880235_2_41:notes"

$notes9 = $slide9.NotesPage
$notes9.Shapes.Item(1).Name = "Google Shape;179;g13a42880235_2_34:notes"
$notes9.Shapes.Item(2).Name = "Google Shape;180;g13a42880235_2_34:notes"
